# Applies the scheduled-runner profit recalculation update across the
# Mateus_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
#
# Each entry below is (SheetName, RowNumber, ColumnNumber, NewValue).
# ColumnNumber maps to the leve-pricing columns H..N (8..14):
#   H=8  currentAveragePrice
#   I=9  currentAveragePriceNQ
#   J=10 currentAveragePriceHQ
#   K=11 LevePriceNQ
#   L=12 LevePriceHQ
#   M=13 LeveProfitNQ
#   N=14 LeveProfitHQ
#
# A New Value of $null means the cell is cleared entirely (removed), matching
# rows where the source data no longer has that column populated.

$wb = $excel.ActiveWorkbook

$updates = @(
    @("ALC", 12, 8, 804.0833),
    @("ALC", 12, 9, 425),
    @("ALC", 12, 10, 2699.5),
    @("ALC", 12, 11, 425),
    @("ALC", 12, 12, 2699.5),
    @("ALC", 12, 13, -255),
    @("ALC", 12, 14, -3039.5),
    @("ALC", 62, 8, 12190.7),
    @("ALC", 62, 9, 14508.714),
    @("ALC", 62, 11, 14508.714),
    @("ALC", 62, 13, -13884.714),
    @("ALC", 64, 8, 9075.9),
    @("ALC", 64, 9, 4151.8),
    @("ALC", 64, 10, 14000),
    @("ALC", 64, 11, 4151.8),
    @("ALC", 64, 12, 14000),
    @("ALC", 64, 13, -3903.8),
    @("ALC", 64, 14, -14496),
    @("ALC", 65, 8, 12190.7),
    @("ALC", 65, 9, 14508.714),
    @("ALC", 65, 11, 72543.57000000001),
    @("ALC", 65, 13, -69423.57000000001),
    @("ALC", 67, 8, 9075.9),
    @("ALC", 67, 9, 4151.8),
    @("ALC", 67, 10, 14000),
    @("ALC", 67, 11, 4151.8),
    @("ALC", 67, 12, 14000),
    @("ALC", 67, 13, -3293.8),
    @("ALC", 67, 14, -15716),
    @("ALC", 80, 8, 1001),
    @("ALC", 80, 9, 1002),
    @("ALC", 80, 10, 1000),
    @("ALC", 80, 11, 3006),
    @("ALC", 80, 12, 3000),
    @("ALC", 80, 13, -2008),
    @("ALC", 80, 14, -4996),
    @("ALC", 83, 8, 1001),
    @("ALC", 83, 9, 1002),
    @("ALC", 83, 10, 1000),
    @("ALC", 83, 11, 9018),
    @("ALC", 83, 12, 9000),
    @("ALC", 83, 13, -4026),
    @("ALC", 83, 14, -18984),
    @("ALC", 107, 8, 2658.5),
    @("ALC", 107, 9, 2658.5),
    @("ALC", 107, 11, 2658.5),
    @("ALC", 107, 13, -738.5),
    @("ALC", 138, 8, 29413870),
    @("ALC", 138, 10, 47620996),
    @("ALC", 138, 12, 142862988),
    @("ALC", 138, 14, -142873268),
    @("ARM", 97, 8, 635.05),
    @("ARM", 97, 9, 774.8333),
    @("ARM", 97, 10, 425.375),
    @("ARM", 97, 11, 774.8333),
    @("ARM", 97, 12, 425.375),
    @("ARM", 97, 13, -278.8333),
    @("ARM", 97, 14, -1417.375),
    @("ARM", 110, 8, 8574.143),
    @("ARM", 110, 9, 7303.9),
    @("ARM", 110, 11, 7303.9),
    @("ARM", 110, 13, -5258.9),
    @("BSM", 94, 8, 2310.4546),
    @("BSM", 94, 9, 2048.8823),
    @("BSM", 94, 11, 2048.8823),
    @("BSM", 94, 13, -1597.8823),
    @("BSM", 134, 8, 2724.625),
    @("BSM", 134, 9, 2706.2666),
    @("BSM", 134, 11, 8118.7998),
    @("BSM", 134, 13, -5583.7998),
    @("CRP", 16, 8, 2758.4092),
    @("CRP", 16, 9, 2128),
    @("CRP", 16, 10, 3861.625),
    @("CRP", 16, 11, 2128),
    @("CRP", 16, 12, 3861.625),
    @("CRP", 16, 13, -1841),
    @("CRP", 16, 14, -4435.625),
    @("CRP", 26, 8, 0),
    @("CRP", 26, 9, 0),
    @("CRP", 26, 11, 0),
    @("CRP", 26, 13, $null),
    @("CRP", 113, 8, 2758.4092),
    @("CRP", 113, 9, 2128),
    @("CRP", 113, 10, 3861.625),
    @("CRP", 113, 11, 2128),
    @("CRP", 113, 12, 3861.625),
    @("CRP", 113, 13, 42),
    @("CRP", 113, 14, -8201.625),
    @("CUL", 12, 8, 999.3570999999999),
    @("CUL", 12, 10, 1144.8),
    @("CUL", 12, 12, 3434.4),
    @("CUL", 12, 14, -3780.4),
    @("CUL", 22, 8, 1000),
    @("CUL", 22, 10, 1000),
    @("CUL", 22, 12, 3000),
    @("CUL", 22, 14, -3338),
    @("CUL", 27, 8, 1000),
    @("CUL", 27, 10, 1000),
    @("CUL", 27, 12, 3000),
    @("CUL", 27, 14, -3204),
    @("CUL", 131, 8, 50002576),
    @("CUL", 131, 10, 5359.75),
    @("CUL", 131, 12, 16079.25),
    @("CUL", 131, 14, -26159.25),
    @("CUL", 132, 8, 38462644),
    @("CUL", 132, 9, 62500980),
    @("CUL", 132, 10, 1309.4),
    @("CUL", 132, 11, 562508820),
    @("CUL", 132, 12, 11784.6),
    @("CUL", 132, 13, -562506290),
    @("CUL", 132, 14, -16844.6),
    @("GSM", 3, 8, 209326.27),
    @("GSM", 3, 10, 37637.5),
    @("GSM", 3, 12, 37637.5),
    @("GSM", 3, 14, -37869.5),
    @("GSM", 7, 8, 4750),
    @("GSM", 7, 9, 4000),
    @("GSM", 7, 10, 5000),
    @("GSM", 7, 11, 4000),
    @("GSM", 7, 12, 5000),
    @("GSM", 7, 13, -3888),
    @("GSM", 7, 14, -5224),
    @("GSM", 8, 8, 4750),
    @("GSM", 8, 9, 4000),
    @("GSM", 8, 10, 5000),
    @("GSM", 8, 11, 4000),
    @("GSM", 8, 12, 5000),
    @("GSM", 8, 13, -3861),
    @("GSM", 8, 14, -5278),
    @("GSM", 9, 8, 2339.6),
    @("GSM", 9, 9, 2299.6667),
    @("GSM", 9, 10, 2399.5),
    @("GSM", 9, 11, 2299.6667),
    @("GSM", 9, 12, 2399.5),
    @("GSM", 9, 13, -2129.6667),
    @("GSM", 9, 14, -2739.5),
    @("GSM", 10, 8, 266.33334),
    @("GSM", 10, 9, 0),
    @("GSM", 10, 10, 266.33334),
    @("GSM", 10, 11, 0),
    @("GSM", 10, 12, 266.33334),
    @("GSM", 10, 13, $null),
    @("GSM", 10, 14, -604.33334),
    @("GSM", 11, 8, 20000000),
    @("GSM", 11, 10, 0),
    @("GSM", 11, 12, 0),
    @("GSM", 11, 14, $null),
    @("GSM", 12, 8, 3999.5),
    @("GSM", 12, 9, 3999.5),
    @("GSM", 12, 11, 3999.5),
    @("GSM", 12, 13, -3859.5),
    @("GSM", 13, 8, 1197),
    @("GSM", 13, 9, 444.8),
    @("GSM", 13, 11, 444.8),
    @("GSM", 13, 13, -305.8),
    @("GSM", 14, 8, 848.9091),
    @("GSM", 14, 9, 733.9),
    @("GSM", 14, 11, 733.9),
    @("GSM", 14, 13, -565.9),
    @("GSM", 20, 8, 16972.334),
    @("GSM", 20, 10, 16972.334),
    @("GSM", 20, 12, 16972.334),
    @("GSM", 20, 14, -17462.334),
    @("GSM", 24, 8, 1512253.5),
    @("GSM", 24, 9, 6000000),
    @("GSM", 24, 10, 16338),
    @("GSM", 24, 11, 6000000),
    @("GSM", 24, 12, 16338),
    @("GSM", 24, 13, -5999827),
    @("GSM", 24, 14, -16684),
    @("GSM", 96, 8, 50395.5),
    @("GSM", 96, 10, 50395.5),
    @("GSM", 96, 12, 50395.5),
    @("GSM", 96, 14, -55887.5),
    @("GSM", 107, 8, 257.5),
    @("GSM", 107, 9, 156.7),
    @("GSM", 107, 10, 425.5),
    @("GSM", 107, 11, 156.7),
    @("GSM", 107, 12, 425.5),
    @("GSM", 107, 13, 1763.3),
    @("GSM", 107, 14, -4265.5),
    @("GSM", 113, 8, 3570.9473),
    @("GSM", 113, 9, 3079.2942),
    @("GSM", 113, 10, 7750),
    @("GSM", 113, 11, 3079.2942),
    @("GSM", 113, 12, 7750),
    @("GSM", 113, 13, -909.2941999999998),
    @("GSM", 113, 14, -12090),
    @("GSM", 139, 8, 0),
    @("GSM", 139, 10, 0),
    @("GSM", 139, 12, 0),
    @("GSM", 139, 14, $null),
    @("LTW", 6, 8, 45000),
    @("LTW", 6, 9, 45000),
    @("LTW", 6, 11, 45000),
    @("LTW", 6, 13, -44888),
    @("LTW", 46, 8, 14465.5625),
    @("LTW", 46, 9, 2492.3333),
    @("LTW", 46, 10, 21649.5),
    @("LTW", 46, 11, 2492.3333),
    @("LTW", 46, 12, 21649.5),
    @("LTW", 46, 13, -2304.3333),
    @("LTW", 46, 14, -22025.5),
    @("LTW", 120, 8, 68993),
    @("LTW", 120, 10, 68993),
    @("LTW", 120, 12, 68993),
    @("LTW", 120, 14, -78669),
    @("LTW", 136, 8, 6534.7),
    @("LTW", 136, 9, 5715.4287),
    @("LTW", 136, 10, 8446.333000000001),
    @("LTW", 136, 11, 17146.2861),
    @("LTW", 136, 12, 25338.999),
    @("LTW", 136, 13, -14596.2861),
    @("LTW", 136, 14, -30438.999),
    @("WVR", 4, 8, 10221.667),
    @("WVR", 4, 9, 4721.6665),
    @("WVR", 4, 11, 4721.6665),
    @("WVR", 4, 13, -4608.6665),
    @("WVR", 14, 8, 15547.774),
    @("WVR", 14, 10, 17845.924),
    @("WVR", 14, 12, 17845.924),
    @("WVR", 14, 14, -18181.924),
    @("WVR", 122, 8, 7100.5557),
    @("WVR", 122, 9, 3968.6667),
    @("WVR", 122, 10, 8666.5),
    @("WVR", 122, 11, 11906.0001),
    @("WVR", 122, 12, 25999.5),
    @("WVR", 122, 13, -9456.000100000001),
    @("WVR", 122, 14, -30899.5)
)

$wsCache = @{}

foreach ($u in $updates) {
    $sheetName = $u[0]
    $row       = [int]$u[1]
    $col       = [int]$u[2]
    $newValue  = $u[3]

    if (-not $wsCache.ContainsKey($sheetName)) {
        $wsCache[$sheetName] = $wb.Worksheets.Item($sheetName)
    }
    $ws = $wsCache[$sheetName]
    $cell = $ws.Cells.Item($row, $col)

    if ($null -eq $newValue) {
        $cell.ClearContents()
    } else {
        $cell.Value = $newValue
    }
}

Write-Host "Applied $($updates.Count) cell updates across $($wsCache.Count) sheets."
